$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are textual (dotted thousands, sub-digit tickers, etc.)
# Force text format so Excel does not silently reinterpret them as numbers/dates.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '42.705.62'
$ws.Range('E2').Value2 = '  +0.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '2.555.80'
$ws.Range('E3').Value2 = '  +0.08%  '
$ws.Range('E4').Value2 = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '302.15'
$ws.Range('E5').Value2 = '  +2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '97.32'
$ws.Range('E6').Value2 = '  +6.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.574'
$ws.Range('E7').Value2 = '  +0.32%  '
$ws.Range('E8').Value2 = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.546'
$ws.Range('E9').Value2 = '  -0.10%  '
$ws.Range('E10').Value2 = '  +2.80%  '
$ws.Range('E11').Value2 = '  +0.38%  '
$ws.Range('E12').Value2 = '  +9.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '7.50'
$ws.Range('E13').Value2 = '  -1.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '2.562.94'
$ws.Range('E14').Value2 = '  +0.92%  '
$ws.Range('E15').Value2 = '  +2.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '14.46'
$ws.Range('E16').Value2 = '  +2.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '42.763.56'
$ws.Range('E17').Value2 = '  +0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '13.49'
$ws.Range('E18').Value2 = '  +8.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '0.0₃0988'
$ws.Range('E19').Value2 = '  +2.02%  '
$ws.Range('E20').Value2 = '  -1.15%  '
$ws.Range('E21').Value2 = '  -1.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '255.84'
$ws.Range('E22').Value2 = '  -0.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '2.94'
$ws.Range('E23').Value2 = '  +2.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '2.08'
$ws.Range('E24').Value2 = '  -1.66%  '
$ws.Range('E25').Value2 = '  -5.59%  '
$ws.Range('E26').Value2 = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '38.99'
$ws.Range('E27').Value2 = '  +8.41%  '
$ws.Range('E28').Value2 = '  +0.50%  '
$ws.Range('E29').Value2 = '  -0.49%  '
$ws.Range('E30').Value2 = '  +1.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '155.83'
$ws.Range('E31').Value2 = '  +3.44%  '
$ws.Range('E32').Value2 = '  -0.20%  '
$ws.Range('E33').Value2 = '  +1.31%  '
$ws.Range('E34').Value2 = '  +1.44%  '
$ws.Range('E35').Value2 = '  -2.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '26.10'
$ws.Range('E36').Value2 = '  +6.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '18.37'
$ws.Range('E37').Value2 = '  +15.52%  '
$ws.Range('E38').Value2 = '  +1.25%  '
$ws.Range('E39').Value2 = '  +0.12%  '
$ws.Range('E40').Value2 = '  +1.26%  '
$ws.Range('E41').Value2 = '  +29.89%  '
$ws.Range('E42').Value2 = '  -2.10%  '
$ws.Range('E43').Value2 = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.998'
$ws.Range('E44').Value2 = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '2.056.40'
$ws.Range('E45').Value2 = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '88.52'
$ws.Range('E46').Value2 = '  +4.97%  '
$ws.Range('E47').Value2 = '  +5.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '76.09'
$ws.Range('E48').Value2 = '  +10.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '2.805.82'
$ws.Range('E49').Value2 = '  +0.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '103.64'
$ws.Range('E50').Value2 = '  +0.40%  '
$ws.Range('E51').Value2 = '  +2.54%  '
